$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column (text) cells
$ws.Range("B2").Value = "<made>"
$ws.Range("B9").Value = "<brop>"
$ws.Range("B13").Value = "<then>"

# Update C column (numeric) cells
$ws.Range("C2").Value = 13
$ws.Range("C3").Value = 13
$ws.Range("C5").Value = 13
$ws.Range("C6").Value = 15
$ws.Range("C7").Value = 14
$ws.Range("C8").Value = 16
$ws.Range("C9").Value = 14
$ws.Range("C10").Value = 13
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 15
$ws.Range("C14").Value = 13
$ws.Range("C15").Value = 14
$ws.Range("C16").Value = 12
$ws.Range("C17").Value = 20
$ws.Range("C18").Value = 12
